$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 36, pushing the former
# rows 36-45 down to 37-46 (dimension grows from A1:R45 to A1:R46).
$ws.Rows.Item(36).Insert()

$ws.Range('A36').Value = 2
$ws.Range('B36').Value = 'Comercializadora del Agro de Limarí'
$ws.Range('C36').Value = 'Coquimbo'
$ws.Range('D36').Value = 44510
$ws.Range('E36').Value = 4
$ws.Range('F36').Value = 100112026
$ws.Range('G36').Value = 'Haba'
$ws.Range('H36').Value = 'Sin especificar'
$ws.Range('I36').Value = 'Primera'
$ws.Range('J36').Value = 1300
$ws.Range('K36').Value = 6000
$ws.Range('L36').Value = 7000
$ws.Range('M36').Value = 6500
$ws.Range('N36').Value = '$/saco 25 kilos'
$ws.Range('O36').Value = 'Provincia de Limarí'
$ws.Range('P36').Value = 260
$ws.Range('Q36').Value = 25
$ws.Range('R36').Value = 'Hortaliza'
